$wb = $excel.ActiveWorkbook

# --- Sheet "P-384": insert a new data row at row 18 (date 42952), shifting subsequent rows down ---
$ws4 = $wb.Worksheets.Item("P-384")

# Extend formatting (date style + border) from the last existing data row down to the new last row
$ws4.Cells.Item(79,1).Copy() | Out-Null
$ws4.Cells.Item(80,1).PasteSpecial(-4122) | Out-Null

# Shift existing rows down by one (processing bottom-up so data is not clobbered before being read)
$ws4.Cells.Item(80,1).Value = 43201
$ws4.Cells.Item(80,2).Value = 0.2928571
$ws4.Cells.Item(79,1).Value = 43175
$ws4.Cells.Item(79,2).Value = 0.403397
$ws4.Cells.Item(78,1).Value = 43166
$ws4.Cells.Item(78,2).Value = 0.6469003
$ws4.Cells.Item(77,1).Value = 43104
$ws4.Cells.Item(77,2).Value = 0.4344392
$ws4.Cells.Item(76,1).Value = 43090
$ws4.Cells.Item(76,2).Value = 0.5082873
$ws4.Cells.Item(75,1).Value = 43073
$ws4.Cells.Item(75,2).Value = 0.5288462
$ws4.Cells.Item(74,1).Value = 43008
$ws4.Cells.Item(74,2).Value = 0.09090910000000001
$ws4.Cells.Item(73,1).Value = 43001
$ws4.Cells.Item(73,2).Value = 0.08883249999999999
$ws4.Cells.Item(72,1).Value = 42992
$ws4.Cells.Item(72,2).Value = 0.09230770000000001
$ws4.Cells.Item(71,1).Value = 42989
$ws4.Cells.Item(71,2).Value = 0.0751445
$ws4.Cells.Item(70,1).Value = 42988
$ws4.Cells.Item(70,2).Value = 0.1021021
$ws4.Cells.Item(69,1).Value = 42979
$ws4.Cells.Item(69,2).Value = 0.1192308
$ws4.Cells.Item(68,1).Value = 42970
$ws4.Cells.Item(68,2).Value = 0.1280488
$ws4.Cells.Item(67,1).Value = 42966
$ws4.Cells.Item(67,2).Value = 0.145749
$ws4.Cells.Item(66,1).Value = 42959
$ws4.Cells.Item(66,2).Value = 0.102439
$ws4.Cells.Item(65,1).Value = 42958
$ws4.Cells.Item(65,2).Value = 0.1542289
$ws4.Cells.Item(64,1).Value = 42951
$ws4.Cells.Item(64,2).Value = 0.1472868
$ws4.Cells.Item(63,1).Value = 42949
$ws4.Cells.Item(63,2).Value = 0.08522730000000001
$ws4.Cells.Item(62,1).Value = 42942
$ws4.Cells.Item(62,2).Value = 0.0763359
$ws4.Cells.Item(61,1).Value = 42941
$ws4.Cells.Item(61,2).Value = 0.1190476
$ws4.Cells.Item(60,1).Value = 42937
$ws4.Cells.Item(60,2).Value = 0.0693069
$ws4.Cells.Item(59,1).Value = 42934
$ws4.Cells.Item(59,2).Value = 0.0743494
$ws4.Cells.Item(58,1).Value = 42926
$ws4.Cells.Item(58,2).Value = 0.0847458
$ws4.Cells.Item(57,1).Value = 42925
$ws4.Cells.Item(57,2).Value = 0.1170213
$ws4.Cells.Item(56,1).Value = 42924
$ws4.Cells.Item(56,2).Value = 0.1273585
$ws4.Cells.Item(55,1).Value = 43225
$ws4.Cells.Item(55,2).Value = 0.1826625
$ws4.Cells.Item(54,1).Value = 43224
$ws4.Cells.Item(54,2).Value = 0.2011494
$ws4.Cells.Item(53,1).Value = 43206
$ws4.Cells.Item(53,2).Value = 0.5033408
$ws4.Cells.Item(52,1).Value = 43192
$ws4.Cells.Item(52,2).Value = 0.4662757
$ws4.Cells.Item(51,1).Value = 43179
$ws4.Cells.Item(51,2).Value = 0.3502415
$ws4.Cells.Item(50,1).Value = 43159
$ws4.Cells.Item(50,2).Value = 0.43898
$ws4.Cells.Item(49,1).Value = 43149
$ws4.Cells.Item(49,2).Value = 0.3625
$ws4.Cells.Item(48,1).Value = 43145
$ws4.Cells.Item(48,2).Value = 0.6065891
$ws4.Cells.Item(47,1).Value = 43138
$ws4.Cells.Item(47,2).Value = 0.4769688
$ws4.Cells.Item(46,1).Value = 43137
$ws4.Cells.Item(46,2).Value = 0.4227941
$ws4.Cells.Item(45,1).Value = 43136
$ws4.Cells.Item(45,2).Value = 0.4263658
$ws4.Cells.Item(44,1).Value = 43095
$ws4.Cells.Item(44,2).Value = 0.3736264
$ws4.Cells.Item(43,1).Value = 43063
$ws4.Cells.Item(43,2).Value = 0.3663194
$ws4.Cells.Item(42,1).Value = 43050
$ws4.Cells.Item(42,2).Value = 0.3011152
$ws4.Cells.Item(41,1).Value = 43048
$ws4.Cells.Item(41,2).Value = 0.2710997
$ws4.Cells.Item(40,1).Value = 43045
$ws4.Cells.Item(40,2).Value = 0.3890411
$ws4.Cells.Item(39,1).Value = 43038
$ws4.Cells.Item(39,2).Value = 0.1458626
$ws4.Cells.Item(38,1).Value = 43037
$ws4.Cells.Item(38,2).Value = 0.1767372
$ws4.Cells.Item(37,1).Value = 43036
$ws4.Cells.Item(37,2).Value = 0.1393324
$ws4.Cells.Item(36,1).Value = 43008
$ws4.Cells.Item(36,2).Value = 0.08053689999999999
$ws4.Cells.Item(35,1).Value = 43007
$ws4.Cells.Item(35,2).Value = 0.1066282
$ws4.Cells.Item(34,1).Value = 43006
$ws4.Cells.Item(34,2).Value = 0.07383969999999999
$ws4.Cells.Item(33,1).Value = 43002
$ws4.Cells.Item(33,2).Value = 0.0841837
$ws4.Cells.Item(32,1).Value = 43001
$ws4.Cells.Item(32,2).Value = 0.08119659999999999
$ws4.Cells.Item(31,1).Value = 43000
$ws4.Cells.Item(31,2).Value = 0.125
$ws4.Cells.Item(30,1).Value = 42992
$ws4.Cells.Item(30,2).Value = 0.125
$ws4.Cells.Item(29,1).Value = 42989
$ws4.Cells.Item(29,2).Value = 0.0826446
$ws4.Cells.Item(28,1).Value = 42986
$ws4.Cells.Item(28,2).Value = 0.07853400000000001
$ws4.Cells.Item(27,1).Value = 42980
$ws4.Cells.Item(27,2).Value = 0.0943396
$ws4.Cells.Item(26,1).Value = 42979
$ws4.Cells.Item(26,2).Value = 0.1176471
$ws4.Cells.Item(25,1).Value = 42975
$ws4.Cells.Item(25,2).Value = 0.1486811
$ws4.Cells.Item(24,1).Value = 42974
$ws4.Cells.Item(24,2).Value = 0.1198157
$ws4.Cells.Item(23,1).Value = 42966
$ws4.Cells.Item(23,2).Value = 0.1318052
$ws4.Cells.Item(22,1).Value = 42957
$ws4.Cells.Item(22,2).Value = 0.1348315
$ws4.Cells.Item(21,1).Value = 42955
$ws4.Cells.Item(21,2).Value = 0.09157510000000001
$ws4.Cells.Item(20,1).Value = 42954
$ws4.Cells.Item(20,2).Value = 0.107438
$ws4.Cells.Item(19,1).Value = 42953
$ws4.Cells.Item(19,2).Value = 0.127551

# Write the newly inserted row
$ws4.Cells.Item(18,1).Value = 42952
$ws4.Cells.Item(18,2).Value = 0.1088083

# --- Sheet "P-891": insert a new data row at row 23 (date 42952), shifting subsequent rows down ---
$ws7 = $wb.Worksheets.Item("P-891")

# Extend formatting (date style + border) from the last existing data row down to the new last row
$ws7.Cells.Item(75,1).Copy() | Out-Null
$ws7.Cells.Item(76,1).PasteSpecial(-4122) | Out-Null

# Shift existing rows down by one (processing bottom-up so data is not clobbered before being read)
$ws7.Cells.Item(76,1).Value = 43189
$ws7.Cells.Item(76,2).Value = 0.4217252
$ws7.Cells.Item(75,1).Value = 43178
$ws7.Cells.Item(75,2).Value = 0.3315927
$ws7.Cells.Item(74,1).Value = 43177
$ws7.Cells.Item(74,2).Value = 0.4324324
$ws7.Cells.Item(73,1).Value = 43176
$ws7.Cells.Item(73,2).Value = 0.4233871
$ws7.Cells.Item(72,1).Value = 43137
$ws7.Cells.Item(72,2).Value = 0.5022971000000001
$ws7.Cells.Item(71,1).Value = 43109
$ws7.Cells.Item(71,2).Value = 0.4695122
$ws7.Cells.Item(70,1).Value = 43073
$ws7.Cells.Item(70,2).Value = 0.4278075
$ws7.Cells.Item(69,1).Value = 43042
$ws7.Cells.Item(69,2).Value = 0.3282548
$ws7.Cells.Item(68,1).Value = 43015
$ws7.Cells.Item(68,2).Value = 0.1403813
$ws7.Cells.Item(67,1).Value = 43014
$ws7.Cells.Item(67,2).Value = 0.1254545
$ws7.Cells.Item(66,1).Value = 43013
$ws7.Cells.Item(66,2).Value = 0.1259542
$ws7.Cells.Item(65,1).Value = 43012
$ws7.Cells.Item(65,2).Value = 0.1945525
$ws7.Cells.Item(64,1).Value = 42993
$ws7.Cells.Item(64,2).Value = 0.1155378
$ws7.Cells.Item(63,1).Value = 42971
$ws7.Cells.Item(63,2).Value = 0.1022727
$ws7.Cells.Item(62,1).Value = 42970
$ws7.Cells.Item(62,2).Value = 0.1158537
$ws7.Cells.Item(61,1).Value = 42961
$ws7.Cells.Item(61,2).Value = 0.1259259
$ws7.Cells.Item(60,1).Value = 42951
$ws7.Cells.Item(60,2).Value = 0.124031
$ws7.Cells.Item(59,1).Value = 42950
$ws7.Cells.Item(59,2).Value = 0.09230770000000001
$ws7.Cells.Item(58,1).Value = 42949
$ws7.Cells.Item(58,2).Value = 0.1420455
$ws7.Cells.Item(57,1).Value = 42942
$ws7.Cells.Item(57,2).Value = 0.0839695
$ws7.Cells.Item(56,1).Value = 42940
$ws7.Cells.Item(56,2).Value = 0.1127451
$ws7.Cells.Item(55,1).Value = 42939
$ws7.Cells.Item(55,2).Value = 0.0695187
$ws7.Cells.Item(54,1).Value = 42937
$ws7.Cells.Item(54,2).Value = 0.0990099
$ws7.Cells.Item(53,1).Value = 42936
$ws7.Cells.Item(53,2).Value = 0.1209964
$ws7.Cells.Item(52,1).Value = 42935
$ws7.Cells.Item(52,2).Value = 0.1046512
$ws7.Cells.Item(51,1).Value = 43261
$ws7.Cells.Item(51,2).Value = 0.0943396
$ws7.Cells.Item(50,1).Value = 43226
$ws7.Cells.Item(50,2).Value = 0.1657609
$ws7.Cells.Item(49,1).Value = 43193
$ws7.Cells.Item(49,2).Value = 0.3198198
$ws7.Cells.Item(48,1).Value = 43189
$ws7.Cells.Item(48,2).Value = 0.3628692
$ws7.Cells.Item(47,1).Value = 43183
$ws7.Cells.Item(47,2).Value = 0.2880658
$ws7.Cells.Item(46,1).Value = 43172
$ws7.Cells.Item(46,2).Value = 0.2525773
$ws7.Cells.Item(45,1).Value = 43164
$ws7.Cells.Item(45,2).Value = 0.3944954
$ws7.Cells.Item(44,1).Value = 43149
$ws7.Cells.Item(44,2).Value = 0.4791667
$ws7.Cells.Item(43,1).Value = 43137
$ws7.Cells.Item(43,2).Value = 0.5514706
$ws7.Cells.Item(42,1).Value = 43104
$ws7.Cells.Item(42,2).Value = 0.5576923
$ws7.Cells.Item(41,1).Value = 43096
$ws7.Cells.Item(41,2).Value = 0.3839662
$ws7.Cells.Item(40,1).Value = 43088
$ws7.Cells.Item(40,2).Value = 0.4149933
$ws7.Cells.Item(39,1).Value = 43066
$ws7.Cells.Item(39,2).Value = 0.2176235
$ws7.Cells.Item(38,1).Value = 43040
$ws7.Cells.Item(38,2).Value = 0.2715827
$ws7.Cells.Item(37,1).Value = 43036
$ws7.Cells.Item(37,2).Value = 0.1872279
$ws7.Cells.Item(36,1).Value = 43018
$ws7.Cells.Item(36,2).Value = 0.1855422
$ws7.Cells.Item(35,1).Value = 43007
$ws7.Cells.Item(35,2).Value = 0.0835735
$ws7.Cells.Item(34,1).Value = 43002
$ws7.Cells.Item(34,2).Value = 0.0841837
$ws7.Cells.Item(33,1).Value = 43001
$ws7.Cells.Item(33,2).Value = 0.0705128
$ws7.Cells.Item(32,1).Value = 43000
$ws7.Cells.Item(32,2).Value = 0.094697
$ws7.Cells.Item(31,1).Value = 42997
$ws7.Cells.Item(31,2).Value = 0.08398949999999999
$ws7.Cells.Item(30,1).Value = 42994
$ws7.Cells.Item(30,2).Value = 0.0932945
$ws7.Cells.Item(29,1).Value = 42993
$ws7.Cells.Item(29,2).Value = 0.0786517
$ws7.Cells.Item(28,1).Value = 42981
$ws7.Cells.Item(28,2).Value = 0.07142859999999999
$ws7.Cells.Item(27,1).Value = 42980
$ws7.Cells.Item(27,2).Value = 0.0825472
$ws7.Cells.Item(26,1).Value = 42979
$ws7.Cells.Item(26,2).Value = 0.1232493
$ws7.Cells.Item(25,1).Value = 42955
$ws7.Cells.Item(25,2).Value = 0.1172161
$ws7.Cells.Item(24,1).Value = 42954
$ws7.Cells.Item(24,2).Value = 0.1487603

# Write the newly inserted row
$ws7.Cells.Item(23,1).Value = 42952
$ws7.Cells.Item(23,2).Value = 0.1243523

$excel.CutCopyMode = $false
